# Refresh the "cryptos" price/volume snapshot (GitHub Actions data pull).
# Every touched cell holds a plain text value (inlineStr in the source
# workbook) -- Coin / Link / Price / Volume(1h) are all formatted as text,
# not numbers, so we assign strings throughout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.803.24"
$ws.Range("E2").Value = "  -0.77%  "
# Row 3
$ws.Range("D3").Value = "1.628.51"
$ws.Range("E3").Value = "  -0.78%  "
# Row 4
$ws.Range("E4").Value = "  +0.04%  "
# Row 5
$ws.Range("D5").Value = "'215.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
# Row 6
$ws.Range("D6").Value = "'0.5071"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "
# Row 7
$ws.Range("E7").Value = "  +0.18%  "
# Row 8
$ws.Range("D8").Value = "'0.2577"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.37%  "
# Row 9
$ws.Range("D9").Value = "'0.06427"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.03%  "
# Row 10
$ws.Range("D10").Value = "'19.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.42%  "
# Row 11
$ws.Range("D11").Value = "'0.07801"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.50%  "
# Row 12
$ws.Range("D12").Value = "'4.258"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.85%  "
# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.631.04"
$ws.Range("E13").Value = "  -0.52%  "
# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.853.35"
$ws.Range("E14").Value = "  -0.78%  "
# Row 15
$ws.Range("D15").Value = "'0.5580"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.70%  "
# Row 16
$ws.Range("D16").Value = "'63.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.89%  "
# Row 17
$ws.Range("D17").Value = "0.0₅7548"
$ws.Range("E17").Value = "  -2.60%  "
# Row 18
$ws.Range("D18").Value = "25.827.48"
# Row 19
$ws.Range("E19").Value = "  +0.06%  "
# Row 20
$ws.Range("D20").Value = "'193.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.44%  "
# Row 21
$ws.Range("D21").Value = "'4.317"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.23%  "
# Row 22
$ws.Range("D22").Value = "'9.817"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.51%  "
# Row 23
$ws.Range("D23").Value = "'6.008"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.96%  "
# Row 24
$ws.Range("E24").Value = "  +0.09%  "
# Row 25
$ws.Range("D25").Value = "'1.793"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.04%  "
# Row 26
$ws.Range("D26").Value = "'141.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.08%  "
# Row 27
$ws.Range("D27").Value = "'0.1264"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "
# Row 28
$ws.Range("D28").Value = "'6.733"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.04%  "
# Row 29
$ws.Range("D29").Value = "'15.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.34%  "
# Row 30
$ws.Range("D30").Value = "'1.237"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.27%  "
# Row 31
$ws.Range("D31").Value = "'0.04873"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.34%  "
# Row 32
$ws.Range("D32").Value = "'3.280"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "
# Row 33
$ws.Range("D33").Value = "'3.191"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.85%  "
# Row 34
$ws.Range("D34").Value = "'1.559"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.31%  "
# Row 35
$ws.Range("D35").Value = "'2.376"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "
# Row 36
$ws.Range("D36").Value = "'0.8955"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.57%  "
# Row 37
$ws.Range("D37").Value = "'2.567"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.25%  "
# Row 38
$ws.Range("D38").Value = "1.129.57"
$ws.Range("E38").Value = "  +2.30%  "
# Row 39
$ws.Range("D39").Value = "'0.5473"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.39%  "
# Row 40
$ws.Range("D40").Value = "'0.01560"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.89%  "
# Row 41
$ws.Range("D41").Value = "'0.9984"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.24%  "
# Row 42
$ws.Range("D42").Value = "'5.575"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.77%  "
# Row 43
$ws.Range("D43").Value = "'0.7966"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.99%  "
# Row 44
$ws.Range("D44").Value = "'97.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.37%  "
# Row 45
$ws.Range("D45").Value = "1.781.85"
$ws.Range("E45").Value = "  +0.17%  "
# Row 46
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  -8.71%  "
# Row 47
$ws.Range("D47").Value = "'0.4433"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.32%  "
# Row 48
$ws.Range("D48").Value = "'55.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.11%  "
# Row 49
$ws.Range("D49").Value = "'0.05054"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.76%  "
# Row 50
$ws.Range("D50").Value = "'7.652"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.08%  "
# Row 51
$ws.Range("D51").Value = "'1.006"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.50%  "
